$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("D3").Value = 44200
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 1400
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = 1450
$ws.Range("P3").Value = 1450

# Row 5 updates
$ws.Range("D5").Value = 44210
$ws.Range("J5").Value = 1450
$ws.Range("K5").Value = 1600
$ws.Range("L5").Value = 1700
$ws.Range("M5").Value = 1650
$ws.Range("P5").Value = 1650
